$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: convert date strings from "YYYY-MM-DD" to "DD/MM/YYYY" text.
# Force Text number format before the write so the COM layer does not
# reinterpret the day/month-ambiguous "DD/MM/YYYY" string as a real date
# serial; restore the style afterwards so no cell formatting changes.
$ws.Range("A276").NumberFormat = "@"
$ws.Range("A276").Value2 = "15/11/2024"
$ws.Range("A276").Style = "Normal"
$ws.Range("A277").NumberFormat = "@"
$ws.Range("A277").Value2 = "02/11/2024"
$ws.Range("A277").Style = "Normal"
$ws.Range("A278").NumberFormat = "@"
$ws.Range("A278").Value2 = "09/11/2024"
$ws.Range("A278").Style = "Normal"
$ws.Range("A279").NumberFormat = "@"
$ws.Range("A279").Value2 = "09/11/2024"
$ws.Range("A279").Style = "Normal"
$ws.Range("A280").NumberFormat = "@"
$ws.Range("A280").Value2 = "17/11/2024"
$ws.Range("A280").Style = "Normal"
$ws.Range("A281").NumberFormat = "@"
$ws.Range("A281").Value2 = "09/11/2024"
$ws.Range("A281").Style = "Normal"
$ws.Range("A282").NumberFormat = "@"
$ws.Range("A282").Value2 = "15/11/2024"
$ws.Range("A282").Style = "Normal"
$ws.Range("A283").NumberFormat = "@"
$ws.Range("A283").Value2 = "15/11/2024"
$ws.Range("A283").Style = "Normal"
$ws.Range("A284").NumberFormat = "@"
$ws.Range("A284").Value2 = "15/11/2024"
$ws.Range("A284").Style = "Normal"
$ws.Range("A285").NumberFormat = "@"
$ws.Range("A285").Value2 = "09/11/2024"
$ws.Range("A285").Style = "Normal"
$ws.Range("A286").NumberFormat = "@"
$ws.Range("A286").Value2 = "15/11/2024"
$ws.Range("A286").Style = "Normal"
$ws.Range("A287").NumberFormat = "@"
$ws.Range("A287").Value2 = "02/11/2024"
$ws.Range("A287").Style = "Normal"
$ws.Range("A288").NumberFormat = "@"
$ws.Range("A288").Value2 = "09/11/2024"
$ws.Range("A288").Style = "Normal"
$ws.Range("A289").NumberFormat = "@"
$ws.Range("A289").Value2 = "06/11/2024"
$ws.Range("A289").Style = "Normal"
$ws.Range("A290").NumberFormat = "@"
$ws.Range("A290").Value2 = "09/11/2024"
$ws.Range("A290").Style = "Normal"
$ws.Range("A291").NumberFormat = "@"
$ws.Range("A291").Value2 = "15/11/2024"
$ws.Range("A291").Style = "Normal"
$ws.Range("A292").NumberFormat = "@"
$ws.Range("A292").Value2 = "09/11/2024"
$ws.Range("A292").Style = "Normal"
$ws.Range("A293").NumberFormat = "@"
$ws.Range("A293").Value2 = "15/11/2024"
$ws.Range("A293").Style = "Normal"
$ws.Range("A294").NumberFormat = "@"
$ws.Range("A294").Value2 = "05/11/2024"
$ws.Range("A294").Style = "Normal"
$ws.Range("A295").NumberFormat = "@"
$ws.Range("A295").Value2 = "05/11/2024"
$ws.Range("A295").Style = "Normal"
$ws.Range("A296").NumberFormat = "@"
$ws.Range("A296").Value2 = "05/10/2024"
$ws.Range("A296").Style = "Normal"
$ws.Range("A297").NumberFormat = "@"
$ws.Range("A297").Value2 = "15/11/2024"
$ws.Range("A297").Style = "Normal"
$ws.Range("A298").NumberFormat = "@"
$ws.Range("A298").Value2 = "15/11/2024"
$ws.Range("A298").Style = "Normal"
$ws.Range("A299").NumberFormat = "@"
$ws.Range("A299").Value2 = "15/11/2024"
$ws.Range("A299").Style = "Normal"
$ws.Range("A300").NumberFormat = "@"
$ws.Range("A300").Value2 = "15/11/2024"
$ws.Range("A300").Style = "Normal"
$ws.Range("A301").NumberFormat = "@"
$ws.Range("A301").Value2 = "02/11/2024"
$ws.Range("A301").Style = "Normal"
$ws.Range("A302").NumberFormat = "@"
$ws.Range("A302").Value2 = "09/11/2024"
$ws.Range("A302").Style = "Normal"
$ws.Range("A303").NumberFormat = "@"
$ws.Range("A303").Value2 = "04/11/2024"
$ws.Range("A303").Style = "Normal"
$ws.Range("A304").NumberFormat = "@"
$ws.Range("A304").Value2 = "08/11/2024"
$ws.Range("A304").Style = "Normal"
$ws.Range("A305").NumberFormat = "@"
$ws.Range("A305").Value2 = "13/11/2024"
$ws.Range("A305").Style = "Normal"
$ws.Range("A306").NumberFormat = "@"
$ws.Range("A306").Value2 = "13/11/2024"
$ws.Range("A306").Style = "Normal"
$ws.Range("A307").NumberFormat = "@"
$ws.Range("A307").Value2 = "13/11/2024"
$ws.Range("A307").Style = "Normal"
$ws.Range("A308").NumberFormat = "@"
$ws.Range("A308").Value2 = "15/11/2024"
$ws.Range("A308").Style = "Normal"
$ws.Range("A309").NumberFormat = "@"
$ws.Range("A309").Value2 = "05/10/2024"
$ws.Range("A309").Style = "Normal"
$ws.Range("A310").NumberFormat = "@"
$ws.Range("A310").Value2 = "15/11/2024"
$ws.Range("A310").Style = "Normal"
$ws.Range("A311").NumberFormat = "@"
$ws.Range("A311").Value2 = "02/11/2024"
$ws.Range("A311").Style = "Normal"
$ws.Range("A312").NumberFormat = "@"
$ws.Range("A312").Value2 = "09/11/2024"
$ws.Range("A312").Style = "Normal"
$ws.Range("A313").NumberFormat = "@"
$ws.Range("A313").Value2 = "15/11/2024"
$ws.Range("A313").Style = "Normal"
$ws.Range("A314").NumberFormat = "@"
$ws.Range("A314").Value2 = "15/11/2024"
$ws.Range("A314").Style = "Normal"
$ws.Range("A315").NumberFormat = "@"
$ws.Range("A315").Value2 = "15/11/2024"
$ws.Range("A315").Style = "Normal"
$ws.Range("A316").NumberFormat = "@"
$ws.Range("A316").Value2 = "15/11/2024"
$ws.Range("A316").Style = "Normal"
$ws.Range("A317").NumberFormat = "@"
$ws.Range("A317").Value2 = "09/11/2024"
$ws.Range("A317").Style = "Normal"
$ws.Range("A318").NumberFormat = "@"
$ws.Range("A318").Value2 = "02/11/2024"
$ws.Range("A318").Style = "Normal"
$ws.Range("A319").NumberFormat = "@"
$ws.Range("A319").Value2 = "09/11/2024"
$ws.Range("A319").Style = "Normal"
$ws.Range("A320").NumberFormat = "@"
$ws.Range("A320").Value2 = "01/11/2024"
$ws.Range("A320").Style = "Normal"
$ws.Range("A321").NumberFormat = "@"
$ws.Range("A321").Value2 = "15/11/2024"
$ws.Range("A321").Style = "Normal"
$ws.Range("A322").NumberFormat = "@"
$ws.Range("A322").Value2 = "02/11/2024"
$ws.Range("A322").Style = "Normal"
$ws.Range("A323").NumberFormat = "@"
$ws.Range("A323").Value2 = "15/11/2024"
$ws.Range("A323").Style = "Normal"
$ws.Range("A324").NumberFormat = "@"
$ws.Range("A324").Value2 = "15/11/2024"
$ws.Range("A324").Style = "Normal"

# Column C: strip the literal " - NA" placeholder segment from Servicio/Producto text
$ws.Range("C276").Value2 = "Fondo - Ahorro"
$ws.Range("C277").Value2 = "Fondo - Almuerzo"
$ws.Range("C278").Value2 = "Fondo - Almuerzo"
$ws.Range("C279").Value2 = "Descuento - Anticipo"
$ws.Range("C280").Value2 = "Descuento - Anticipo"
$ws.Range("C283").Value2 = "Fondo - Ahorro"
$ws.Range("C284").Value2 = "Apoyo - 11 días"
$ws.Range("C285").Value2 = "Fondo - Almuerzo"
$ws.Range("C286").Value2 = "Fondo - Ahorro"
$ws.Range("C287").Value2 = "Fondo - Almuerzo"
$ws.Range("C288").Value2 = "Fondo - Almuerzo"
$ws.Range("C289").Value2 = "Descuento - Anticipo"
$ws.Range("C290").Value2 = "Descuento - Anticipo"
$ws.Range("C291").Value2 = "Fondo - Ahorro"
$ws.Range("C292").Value2 = "Fondo - Almuerzo"
$ws.Range("C297").Value2 = "Fondo - Prestamo"
$ws.Range("C298").Value2 = "Fondo - Intereses"
$ws.Range("C299").Value2 = "Fondo - Prestamo"
$ws.Range("C300").Value2 = "Fondo - Intereses"
$ws.Range("C301").Value2 = "Fondo - Almuerzo"
$ws.Range("C302").Value2 = "Fondo - Almuerzo"
$ws.Range("C303").Value2 = "Descuento - Anticipo"
$ws.Range("C304").Value2 = "Descuento - Anticipo"
$ws.Range("C305").Value2 = "Descuento - Anticipo"
$ws.Range("C306").Value2 = "Descuento - Anticipo"
$ws.Range("C308").Value2 = "Fondo - Ahorro"
$ws.Range("C310").Value2 = "Fondo - Ahorro"
$ws.Range("C311").Value2 = "Fondo - Almuerzo"
$ws.Range("C312").Value2 = "Fondo - Almuerzo"
$ws.Range("C313").Value2 = "Fondo - Ahorro"
$ws.Range("C314").Value2 = "Fondo - Ahorro"
$ws.Range("C315").Value2 = "Fondo - Prestamo"
$ws.Range("C316").Value2 = "Fondo - Intereses"
$ws.Range("C317").Value2 = "Fondo - Almuerzo"
$ws.Range("C318").Value2 = "Descuento - Anticipo"
$ws.Range("C319").Value2 = "Fondo - desayuno"
$ws.Range("C321").Value2 = "Fondo - Ahorro"
$ws.Range("C322").Value2 = "Fondo - Almuerzo"
$ws.Range("C324").Value2 = "Fondo - Intereses"
